$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 15 (row 19): record the Part I solve time.
$ws.Range("D19").Value = 3

# Leave the selection where Excel would land after typing the value and
# pressing Enter (one row down, same column).
$ws.Range("D20").Select() | Out-Null
